$d = $word.ActiveDocument

# ============================================================
# EDIT 1: Merge the "Languages: ..." and "Tools: ..." bullet
#         paragraphs (under "Technical Skills") into a single
#         "Languages and Tools: ..." paragraph. Also trims the
#         tools list (drops "MySQL", "Handlebars", "Query") and
#         appends ", Agile Methodology".
# ============================================================

$langIndex = -1
$count = $d.Paragraphs.Count
for ($i = 1; $i -le $count; $i++) {
    $t = $d.Paragraphs($i).Range.Text
    if ($t.StartsWith("Languages:")) {
        $langIndex = $i
        break
    }
}
if ($langIndex -eq -1) {
    throw "Could not find the 'Languages:' paragraph"
}

$pLang = $d.Paragraphs($langIndex)
$pTools = $d.Paragraphs($langIndex + 1)
if (-not $pTools.Range.Text.StartsWith("Tools:")) {
    throw "Paragraph after 'Languages:' was not the expected 'Tools:' paragraph"
}

# Merge the two paragraphs by deleting the paragraph mark that ends the
# "Languages" paragraph (this is what Word does when you press Delete at
# the end of a paragraph to join it with the next one).
$markRange = $d.Range($pLang.Range.End - 1, $pLang.Range.End)
$markRange.Delete()

# Re-fetch the merged paragraph and locate the old "Tools: ..." text inside
# it so the edits below are anchored to real offsets rather than guesses.
$pMerged = $d.Paragraphs($langIndex)
$base = $pMerged.Range.Start
$mergedText = $pMerged.Range.Text

$oldTools = "Tools: GitHub, MongoDB, MySQL, Express, React, Node, Handlebars, Query, Bootstrap, SDLC"
$toolsStart = $mergedText.IndexOf($oldTools)
if ($toolsStart -lt 0) {
    throw "Could not find the expected tools text in the merged paragraph"
}

$labelLen = "Tools: ".Length
$skillsListStart = $base + $toolsStart + $labelLen
$skillsListEnd = $base + $toolsStart + $oldTools.Length

# 1) Replace the tools list text (keeps the original run's formatting -
#    Roboto, matching the rest of the skills list).
$rSkills = $d.Range($skillsListStart, $skillsListEnd)
$rSkills.Text = "GitHub, MongoDB, Express, React, Node, Bootstrap, SDLC, Agile Methodology"

# 2) Remove the now-redundant "Tools: " label that used to start the
#    second bullet.
$rLabel = $d.Range($base + $toolsStart, $skillsListStart)
$rLabel.Text = ""

# 3) Insert the ", " separator between the languages list and the tools
#    list so they read as one comma-separated run.
$sepPos = $base + $toolsStart
$rSep = $d.Range($sepPos, $sepPos)
$rSep.InsertBefore(", ")

# 4) Update the "Languages: " heading to "Languages and Tools: ".
$rLangLabel = $d.Range($base, $base + "Languages: ".Length)
$rLangLabel.Text = "Languages and Tools: "

Write-Host "Languages/Tools paragraph now reads: $($d.Paragraphs($langIndex).Range.Text)"

# ============================================================
# EDIT 2: Rewrite the "A 24-week intensive program..." bullet
#         under the coding bootcamp entry.
# ============================================================

$progIndex = -1
$count = $d.Paragraphs.Count
for ($i = 1; $i -le $count; $i++) {
    $t = $d.Paragraphs($i).Range.Text
    if ($t.StartsWith("A 24-week intensive program")) {
        $progIndex = $i
        break
    }
}
if ($progIndex -eq -1) {
    throw "Could not find the 'A 24-week intensive program...' paragraph"
}

$pProg = $d.Paragraphs($progIndex)
$progFull = $pProg.Range.Text
$progNoMark = $progFull.Substring(0, $progFull.Length - 1)
$rProg = $d.Range($pProg.Range.Start, $pProg.Range.Start + $progNoMark.Length)
$rProg.Text = "24-week intensive program focused on gaining technical programming skills including HTML5, CSS3, JavaScript, jQuery, Progressive Web Apps, Agile Methodology, Bootstrap, React.js, Express.js, Node.js, MongoDB, MySQL, Command Line, OOP, Git, Python, Java, C#, Amazon Web Services and more"

Write-Host "Program paragraph now reads: $($d.Paragraphs($progIndex).Range.Text)"
